# TOR-57 Add handle multi die
# Refactor the "Reaction Roll" sheet's 2d6 lookup table so a single roll
# column (sum of the dice) maps straight to a result, with bucketed ranges
# shown for the ambiguous outcomes, instead of enumerating every individual
# (die1, die2) combination.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reaction Roll")

# Drop the old "roll1_d6" / "roll2_d6" columns (C & D) - only the roll sum
# and the result text are kept going forward.
$ws.Range("C1:D1").EntireColumn.Delete()

# Collapse the 36-row combination table down to the 6 rows we still need.
$ws.Range("A7:A37").EntireRow.Delete()

# Row 1 - headers
$ws.Range("A1").Value = "2d6"
$ws.Range("B1").Value = "Result"

# Row 2 - Hostile
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Hostile"

# Row 3 - Unfriendly
$ws.Range("A3").Value = "[3,4,5]"
$ws.Range("B3").Value = "Unfriendly"

# Row 4 - Unsure
$ws.Range("A4").Value = "[6,7,8]"
$ws.Range("B4").Value = "Unsure. "

# Row 5 - Talkative
$ws.Range("A5").Value = "[9,10,11]"
$ws.Range("B5").Value = "Talkative. "

# Row 6 - Helpful
$ws.Range("A6").Value = 12
$ws.Range("B6").Value = "Helpful. "

# Match the author's recorded selection/cursor position after the edit.
$null = $ws.Range("A6").Select()
